# Update odds values on Sheet1 to reflect the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (Kuressaare - Flora)
$ws.Range("G5").Value = 6.8
$ws.Range("H5").Value = 4.85
$ws.Range("I5").Value = 1.32
$ws.Range("V5").Value = 18
$ws.Range("W5").Value = 110
$ws.Range("X5").Value = 50
$ws.Range("Y5").Value = 40
$ws.Range("AA5").Value = 8.75
$ws.Range("AB5").Value = 14.5
$ws.Range("AD5").Value = 250
$ws.Range("AF5").Value = 6.6
$ws.Range("AH5").Value = 7.8

# Row 6 (Tammeka - Tallinna Kalev)
$ws.Range("G6").Value = 2.25
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.7
$ws.Range("T6").Value = 10.25
$ws.Range("U6").Value = 12.5
$ws.Range("V6").Value = 7.8
$ws.Range("W6").Value = 21
$ws.Range("AA6").Value = 6.4
$ws.Range("AB6").Value = 9
$ws.Range("AE6").Value = 10.75
$ws.Range("AF6").Value = 14
$ws.Range("AG6").Value = 8.5
$ws.Range("AH6").Value = 27

# Row 8 (Braunschweig - Saarbrucken)
$ws.Range("G8").Value = 1.83
$ws.Range("I8").Value = 4.1
$ws.Range("U8").Value = 9
$ws.Range("AD8").Value = 251
$ws.Range("AE8").Value = 12
$ws.Range("AF8").Value = 21
$ws.Range("AH8").Value = 41

# Row 11 (Suwon FC - Jeju SK)
$ws.Range("G11").Value = 2.55
$ws.Range("I11").Value = 2.75
$ws.Range("K11").Value = 9.5
$ws.Range("T11").Value = 8.5
$ws.Range("U11").Value = 12
$ws.Range("Y11").Value = 29
$ws.Range("AA11").Value = 6.5
$ws.Range("AH11").Value = 29
$ws.Range("AI11").Value = 23

# Row 12 (Sandviken - Sundsvall)
$ws.Range("N12").Value = 1.83
$ws.Range("O12").Value = 1.98

# Row 14 (Grasshoppers - Aarau)
$ws.Range("G14").Value = 1.75
$ws.Range("H14").Value = 3.75
$ws.Range("I14").Value = 4.2
